# TAOrderExample.xlsx update:
#  - the "TagsConfig" row no longer ships a fog override, so clear B6
#    (this also drops the now-unused shared string on save)
#  - widen column B so the longer config paths are readable
#  - leave the sheet scrolled/zoomed in on the TagsConfig row, matching
#    where the author was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "fog=true/fog_visible=85" tag value from B6 (TagsConfig row)
$ws.Range("B6").ClearContents()

# Widen column B (was ~30.63 chars, now ~52.63 chars wide)
$ws.Columns.Item(2).ColumnWidth = 51.83

# Leave the view zoomed to 115% with B6 selected/active
$ws.Range("B6").Select()
$excel.ActiveWindow.Zoom = 115
